$d = $word.ActiveDocument

function New-CollapsedRange([int]$pos) {
    return $d.Range($pos, $pos)
}

function Insert-BodyXmlAt($range, [string]$bodyXml) {
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# 1) Insert six new paragraphs right after "Tst\SampleProtocols.txt" and
#    before the two blank paragraphs that follow it.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$found = $anchor.Find.Execute("Tst\SampleProtocols.txt", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "anchor paragraph not found" }
$anchorPara = $anchor.Paragraphs(1)
$anchorPara.Range.InsertParagraphAfter()
$newFirstPara = $anchorPara.Next()
$insertPoint = New-CollapsedRange $newFirstPara.Range.Start

$newParasXml = @'
<w:p><w:r><w:t>Please note that in the file with test directories, the names of the directories should either be absolute paths, or paths relative to Tst directory. IN the latter case, the directory name should start with either “.\” or with the name (no backslash). For example:</w:t></w:r><w:r><w:br/><w:t>Correct directory names:</w:t></w:r></w:p><w:p><w:r><w:t>D:\PLanguage\plang\Tst\RegressionTests</w:t></w:r></w:p><w:p><w:r><w:t>.</w:t></w:r><w:r><w:t>\RegressionTests</w:t></w:r></w:p><w:p><w:r><w:t>RegressionTests</w:t></w:r></w:p><w:p><w:r><w:t>Incorrect directory name:</w:t></w:r></w:p><w:p><w:r><w:t>\</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>RegressionTests</w:t></w:r></w:p>
'@

Insert-BodyXmlAt $insertPoint $newParasXml

Write-Output "Step 1 done: new paragraphs inserted"
